$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column I (9), shifting existing columns I..O to K..Q
$ws.Range("I1:J1").EntireColumn.Insert()

# Set header values for the two new columns
$ws.Range("I1").Value = "KONTAKTNÍ_OSOBA"
$ws.Range("J1").Value = "KONTAKTNÍ_INFO"

# Set column widths for the two new columns (target stored widths are
# 25.28515625 / 28.7109375; ColumnWidth is specified in character units,
# which is stored width minus ~5/6 of a character)
$ws.Range("I1").EntireColumn.ColumnWidth = 24.451822916666668
$ws.Range("J1").EntireColumn.ColumnWidth = 27.877604166666668

# Set selection to match target state
$ws.Range("A2:XFD1048576").Select()
